# NIT-9007995718 Estado de Cuenta: remove previous account-statement rows and
# regenerate the worker/period table with the latest data; update the summary
# totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: the table body grows from 19 data rows (16-34) to 22 data
#       rows (16-37). Insert 3 rows right after the current last data row (34)
#       so the trailing gap + signature block shift down automatically.
$ws.Rows("35:37").Insert()

# --- 2. Formatting, done BEFORE overwriting values: row 34 still carries the
#       special "last row of the table" border (it hasn't moved yet) - copy
#       that look onto the new final row (37, Oscar's record after the
#       regeneration) first, then stamp the regular interior-row look (still
#       on row 16) across the whole body so row 34 (no longer last) reverts
#       to the normal look too.
$ws.Range("B34:J34").Copy()
$ws.Range("B37:J37").PasteSpecial(-4122) | Out-Null

$ws.Range("B16:J16").Copy()
$ws.Range("B16:J36").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- 3. New full table body (rows 16-37): Tipo Doc, N Doc, Nombre, Periodo,
#       Valor Mora, Salario Basico
$data = @(
    ,@(16, 'CC', '73190731', 'CARLOS ENRIQUE TAFUR RODRIGUEZ', '2501', 56940, 877803)
    ,@(17, 'CC', '73190731', 'CARLOS ENRIQUE TAFUR RODRIGUEZ', '2412', 56940, 877803)
    ,@(18, 'CC', '91523749', 'JUAN CARLOS REY MARTINEZ', '2304', 46400, 1600000)
    ,@(19, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2507', 56000, 1400000)
    ,@(20, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2506', 56000, 1400000)
    ,@(21, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2505', 56000, 1400000)
    ,@(22, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2504', 56000, 1400000)
    ,@(23, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2503', 56000, 1400000)
    ,@(24, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2502', 56000, 1400000)
    ,@(25, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2501', 56000, 1400000)
    ,@(26, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2412', 56000, 1400000)
    ,@(27, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2411', 56000, 1400000)
    ,@(28, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2410', 56000, 1400000)
    ,@(29, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2409', 56000, 1400000)
    ,@(30, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2408', 56000, 1400000)
    ,@(31, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2407', 56000, 1400000)
    ,@(32, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2406', 56000, 1400000)
    ,@(33, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2405', 56000, 1400000)
    ,@(34, 'CC', '84453346', 'STEVE ALBERTO GARCIA REVOLLEDO', '2404', 56000, 1400000)
    ,@(35, 'CC', '1002202553', 'MARIA JOSE CASTRO RAMIREZ', '2403', 56000, 1400000)
    ,@(36, 'CC', '1002202553', 'MARIA JOSE CASTRO RAMIREZ', '2402', 18667, 1400000)
    ,@(37, 'CC', '1051363446', 'OSCAR DAVID JULIO ELKAIEK', '2506', 28470, 1423500)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("F$r").Value = $row[5]
    $ws.Range("G$r").Value = $row[6]
}

# --- 4. Header / summary cells
$ws.Range("E11").Value = 1159417      # VALOR MORA total
$ws.Range("C13").Value = 5            # Cant. Trabajadores
$ws.Range("F13").Value = 19           # Cant. Periodos
